$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) from the last existing row (226) down into the
# three new rows, then overwrite the values with the updated figures.
$ws.Range("A226:D226").Copy($ws.Range("A227:D227"))
$ws.Range("A226:D226").Copy($ws.Range("A228:D228"))
$ws.Range("A226:D226").Copy($ws.Range("A229:D229"))

$ws.Range("A227").Value = 44301
$ws.Range("B227").Value = 2
$ws.Range("C227").Value = 2
$ws.Range("D227").Value = 107.469102632993

$ws.Range("A228").Value = 44302
$ws.Range("B228").Value = 0
$ws.Range("C228").Value = 2
$ws.Range("D228").Value = 107.469102632993

$ws.Range("A229").Value = 44303
$ws.Range("B229").Value = 0
$ws.Range("C229").Value = 2
$ws.Range("D229").Value = 107.469102632993
